$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New column headers (set D1 before C1 so shared-string indices match target order)
$ws.Range("D1").Value = "截至20180722的文件数"
$ws.Range("C1").Value = "footer"

# New column widths (closest achievable values given engine's column-width quantization)
$ws.Range("C1").EntireColumn.ColumnWidth = 17.145
$ws.Range("D1").EntireColumn.ColumnWidth = 19.715

# New data values
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 592

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 395

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 74

# Update selection to match target state
$ws.Range("D15").Select()
